# Updated symbol list on Wed Dec 28 22:26:06 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numeric-looking values as plain TEXT in the
# source workbook (no leading apostrophe shown in the UI, just a text cell).
# A bare `Range.Value = "243.12"` assignment would let Excel's type-inference
# re-interpret that text as a Number, which would NOT match the original
# file's cell type. To keep the cell a text cell (same as before the edit),
# we briefly mark the cell as Text-formatted before writing the value, then
# clear the formatting again so the cell's style returns to the workbook's
# default (unstyled) look - only the stored value/type actually changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.ClearFormats()
}

# --- Column D ("Price") updates ---
Set-TextValue "D2"  "243.12"
Set-TextValue "D3"  "23.74"
Set-TextValue "D4"  "5.232"
Set-TextValue "D5"  "0.05758"
Set-TextValue "D6"  "6.407"
Set-TextValue "D7"  "3.229"
Set-TextValue "D8"  "0.8062"
Set-TextValue "D9"  "0.8798"
Set-TextValue "D10" "0.1374"
Set-TextValue "D11" "0.07075"
Set-TextValue "D12" "0.03141"
Set-TextValue "D13" "0.03045"
Set-TextValue "D14" "0.09326"
Set-TextValue "D15" "3.819"
Set-TextValue "D16" "0.001519"
Set-TextValue "D17" "0.04699"
Set-TextValue "D18" "0.0006008"
Set-TextValue "D19" "0.006258"
Set-TextValue "D20" "0.001260"
Set-TextValue "D21" "0.004053"
Set-TextValue "D22" "0.00008714"
Set-TextValue "D23" "3.547"
Set-TextValue "D24" "2.154"
Set-TextValue "D25" "0.3161"
Set-TextValue "D28" "0.0002332"
Set-TextValue "D40" "0.03727"
Set-TextValue "D41" "0.006250"
Set-TextValue "D42" "0.1045"
Set-TextValue "D43" "0.002501"
Set-TextValue "D44" "0.007164"
Set-TextValue "D45" "0.00005332"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "D47" "0.5358"
Set-TextValue "D48" "0.002438"
Set-TextValue "D49" "0.00002103"
Set-TextValue "D50" "0.0002003"

# --- Column E ("Volume(1h)") updates ---
# These are non-numeric strings, so plain assignment keeps them as text.
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("E48").Value = "47BOLOBOLO"

Write-Host "Applied cryptos.xlsx symbol-list update"
